$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and volume-change percentage (column E)
# values with the latest scraped figures. A leading apostrophe is used on the
# Price column assignments so that values which look numeric (e.g. "316.64")
# are stored as text, matching the original inline-string cell type.
$ws.Range("D2").Value = '''24.735.78'
$ws.Range("E2").Value = '  +1.63%  '
$ws.Range("D3").Value = '''1.696.03'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = '''316.64'
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").Value = '''0.3948'
$ws.Range("E7").Value = '  +0.59%  '
$ws.Range("D8").Value = '''0.4058'
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("D9").Value = '''1.486'
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("D10").Value = '''1.000'
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("D11").Value = '''53.38'
$ws.Range("E11").Value = '  -1.81%  '
$ws.Range("D12").Value = '''0.08875'
$ws.Range("E12").Value = '  +1.45%  '
$ws.Range("D13").Value = '''7.234'
$ws.Range("E13").Value = '  -1.07%  '
$ws.Range("D14").Value = '''23.65'
$ws.Range("E14").Value = '  +2.43%  '
$ws.Range("D15").Value = '''8.047'
$ws.Range("E15").Value = '  +8.58%  '
$ws.Range("D16").Value = '''0.00001322'
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").Value = '''1.693.61'
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("D18").Value = '''100.02'
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").Value = '''0.07013'
$ws.Range("E19").Value = '  -0.43%  '
$ws.Range("E20").Value = '  +1.49%  '
$ws.Range("D21").Value = '''7.072'
$ws.Range("E21").Value = '  +4.95%  '
$ws.Range("D22").Value = '''1.000'
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("E23").Value = '  +1.41%  '
$ws.Range("D24").Value = '''24.733.25'
$ws.Range("E24").Value = '  +1.64%  '
$ws.Range("D25").Value = '''3.260'
$ws.Range("E25").Value = '  +9.59%  '
$ws.Range("D26").Value = '''2.352'
$ws.Range("E26").Value = '  +1.69%  '
$ws.Range("D27").Value = '''22.75'
$ws.Range("E27").Value = '  +2.09%  '
$ws.Range("D28").Value = '''163.23'
$ws.Range("E28").Value = '  +2.27%  '
$ws.Range("D29").Value = '''136.26'
$ws.Range("E29").Value = '  +2.13%  '
$ws.Range("E30").Value = '  +1.76%  '
$ws.Range("D31").Value = '''7.487'
$ws.Range("E31").Value = '  +1.26%  '
$ws.Range("D32").Value = '''1.879.65'
$ws.Range("E32").Value = '  -0.54%  '
$ws.Range("D33").Value = '''1.070'
$ws.Range("E33").Value = '  -0.81%  '
$ws.Range("D34").Value = '''0.08595'
$ws.Range("E34").Value = '  -1.25%  '
$ws.Range("D35").Value = '''7.151'
$ws.Range("E35").Value = '  -3.26%  '
$ws.Range("D36").Value = '''11.61'
$ws.Range("E36").Value = '  +5.47%  '
$ws.Range("E37").Value = '  +2.25%  '
$ws.Range("D38").Value = '''1.931'
$ws.Range("E38").Value = '  -0.79%  '
$ws.Range("D39").Value = '''14.51'
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("D40").Value = '''0.09239'
$ws.Range("E40").Value = '  +3.40%  '
$ws.Range("D41").Value = '''0.02738'
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("D42").Value = '''1.473'
$ws.Range("E42").Value = '  +1.02%  '
$ws.Range("D43").Value = '''0.7682'
$ws.Range("E43").Value = '  +1.54%  '
$ws.Range("D44").Value = '''15.94'
$ws.Range("E44").Value = '  +3.43%  '
$ws.Range("D45").Value = '''0.7204'
$ws.Range("E45").Value = '  +1.17%  '
$ws.Range("D46").Value = '''2.578'
$ws.Range("E46").Value = '  +5.86%  '
$ws.Range("D47").Value = '''4.224'
$ws.Range("E47").Value = '  +2.29%  '
$ws.Range("D48").Value = '''0.9997'
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("D49").Value = '''1.326'
$ws.Range("E49").Value = '  +3.55%  '
$ws.Range("D50").Value = '''139.52'
$ws.Range("E50").Value = '  -0.32%  '
$ws.Range("D51").Value = '''0.07989'
$ws.Range("E51").Value = '  +0.78%  '
